$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-6
# from 45175 (2023-09-06) to 45183 (2023-09-14)
$newDate = [DateTime]::FromOADate(45183)

$ws.Range("C2:C6").Value = $newDate
